$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Dauer" -> "Dauer in h"
$ws.Range("D2").Value = "Dauer in h"

# New total row: sum of all duration entries
$ws.Range("D3").Formula = "=SUM(D4:D48)"

# Convert duration entries from text ("Xh" / "Xh Ymin") to plain numeric hour values
$ws.Range("D4").Value = 7
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 4
$ws.Range("D8").Value = 2
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("D21").Value = 2
$ws.Range("D23").Value = 8
$ws.Range("D25").Value = 10
$ws.Range("D28").Value = 4
$ws.Range("D29").Value = 2
$ws.Range("D30").Value = 3
$ws.Range("D32").Value = 1
$ws.Range("D33").Value = 2
$ws.Range("D34").Value = 1
$ws.Range("D36").Value = 4
$ws.Range("D37").Value = 4
$ws.Range("D39").Value = 5
$ws.Range("D41").Value = 5
$ws.Range("D42").Value = 2
$ws.Range("D44").Value = 3
$ws.Range("D46").Value = 2
$ws.Range("D48").Value = 4

# Update the view selection (also resets the scrolled-down top-left cell)
$ws.Range("D2").Select()
